$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 6118
$ws.Cells.Item(43, 9).Value = 6865.6
$ws.Cells.Item(43, 10).Value = 4249
$ws.Cells.Item(43, 11).Value = 6865.6
$ws.Cells.Item(43, 12).Value = 4249
$ws.Cells.Item(43, 13).Value = -6796.6
$ws.Cells.Item(43, 14).Value = -4387
$ws.Cells.Item(101, 8).Value = 1389.25
$ws.Cells.Item(101, 9).Value = 357
$ws.Cells.Item(101, 11).Value = 1071
$ws.Cells.Item(101, 13).Value = 551
$ws.Cells.Item(111, 8).Value = 2118.6
$ws.Cells.Item(111, 9).Value = 1050.5
$ws.Cells.Item(111, 10).Value = 2830.6667
$ws.Cells.Item(111, 11).Value = 3151.5
$ws.Cells.Item(111, 12).Value = 8492.000100000001
$ws.Cells.Item(111, 13).Value = -84.5
$ws.Cells.Item(111, 14).Value = -14626.0001
$ws.Cells.Item(131, 8).Value = 418334.84
$ws.Cells.Item(131, 9).Value = 418334.84
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 1255004.52
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 13).Value = -1249964.52
$ws.Cells.Item(131, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 1797.5657
$ws.Cells.Item(138, 9).Value = 1187.2413
$ws.Cells.Item(138, 10).Value = 2050.4143
$ws.Cells.Item(138, 11).Value = 3561.7239
$ws.Cells.Item(138, 12).Value = 6151.242899999999
$ws.Cells.Item(138, 13).Value = 1578.2761
$ws.Cells.Item(138, 14).Value = -16431.2429

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2418.9102
$ws.Cells.Item(32, 9).Value = 1827.6173
$ws.Cells.Item(32, 10).Value = 8405.75
$ws.Cells.Item(32, 11).Value = 1827.6173
$ws.Cells.Item(32, 12).Value = 8405.75
$ws.Cells.Item(32, 13).Value = -1540.6173
$ws.Cells.Item(32, 14).Value = -8979.75
$ws.Cells.Item(92, 8).Value = 40000
$ws.Cells.Item(92, 10).Value = 40000
$ws.Cells.Item(92, 12).Value = 40000
$ws.Cells.Item(92, 14).Value = -44992
$ws.Cells.Item(102, 8).Value = 2921.3635
$ws.Cells.Item(102, 9).Value = 3013.5
$ws.Cells.Item(102, 11).Value = 3013.5
$ws.Cells.Item(102, 13).Value = -1391.5
$ws.Cells.Item(110, 8).Value = 5223.923
$ws.Cells.Item(110, 9).Value = 6745.4443
$ws.Cells.Item(110, 11).Value = 6745.4443
$ws.Cells.Item(110, 13).Value = -4700.4443
$ws.Cells.Item(132, 8).Value = 20824.576
$ws.Cells.Item(132, 9).Value = 6685.9
$ws.Cells.Item(132, 10).Value = 29661.25
$ws.Cells.Item(132, 11).Value = 20057.7
$ws.Cells.Item(132, 12).Value = 88983.75
$ws.Cells.Item(132, 13).Value = -17527.7
$ws.Cells.Item(132, 14).Value = -94043.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 31462.814
$ws.Cells.Item(20, 9).Value = 11453
$ws.Cells.Item(20, 10).Value = 43233.293
$ws.Cells.Item(20, 11).Value = 11453
$ws.Cells.Item(20, 12).Value = 43233.293
$ws.Cells.Item(20, 13).Value = -11206
$ws.Cells.Item(20, 14).Value = -43727.293
$ws.Cells.Item(94, 8).Value = 1089.1936
$ws.Cells.Item(94, 9).Value = 629.75
$ws.Cells.Item(94, 11).Value = 629.75
$ws.Cells.Item(94, 13).Value = -178.75
$ws.Cells.Item(105, 8).Value = 1980
$ws.Cells.Item(105, 9).Value = 1990
$ws.Cells.Item(105, 10).Value = 1920
$ws.Cells.Item(105, 11).Value = 1990
$ws.Cells.Item(105, 12).Value = 1920
$ws.Cells.Item(105, 13).Value = -243
$ws.Cells.Item(105, 14).Value = -5414

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1200.6666
$ws.Cells.Item(22, 10).Value = 2014.25
$ws.Cells.Item(22, 12).Value = 2014.25
$ws.Cells.Item(22, 14).Value = -2714.25
$ws.Cells.Item(58, 8).Value = 15275.3955
$ws.Cells.Item(58, 9).Value = 5496.6562
$ws.Cells.Item(58, 11).Value = 5496.6562
$ws.Cells.Item(58, 13).Value = -5293.6562
$ws.Cells.Item(105, 8).Value = 51005
$ws.Cells.Item(132, 8).Value = 6825.0415
$ws.Cells.Item(132, 9).Value = 2489.1904
$ws.Cells.Item(132, 11).Value = 7467.5712
$ws.Cells.Item(132, 13).Value = -4937.5712
$ws.Cells.Item(136, 8).Value = 15275.3955
$ws.Cells.Item(136, 9).Value = 5496.6562
$ws.Cells.Item(136, 11).Value = 16489.9686
$ws.Cells.Item(136, 13).Value = -13939.9686

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 50.238094
$ws.Cells.Item(12, 10).Value = 45.583332
$ws.Cells.Item(12, 12).Value = 136.749996
$ws.Cells.Item(12, 14).Value = -482.749996
$ws.Cells.Item(63, 8).Value = 40403
$ws.Cells.Item(63, 10).Value = 40403
$ws.Cells.Item(63, 12).Value = 121209
$ws.Cells.Item(63, 14).Value = -122707
$ws.Cells.Item(66, 8).Value = 40403
$ws.Cells.Item(66, 10).Value = 40403
$ws.Cells.Item(66, 12).Value = 363627
$ws.Cells.Item(66, 14).Value = -371115
$ws.Cells.Item(68, 8).Value = 1728.0714
$ws.Cells.Item(68, 10).Value = 2160.1
$ws.Cells.Item(68, 12).Value = 6480.299999999999
$ws.Cells.Item(68, 14).Value = -8102.299999999999
$ws.Cells.Item(71, 8).Value = 1728.0714
$ws.Cells.Item(71, 10).Value = 2160.1
$ws.Cells.Item(71, 12).Value = 19440.9
$ws.Cells.Item(71, 14).Value = -27552.9
$ws.Cells.Item(113, 8).Value = 1013.2593
$ws.Cells.Item(113, 10).Value = 1000.5
$ws.Cells.Item(113, 12).Value = 3001.5
$ws.Cells.Item(113, 14).Value = -7341.5
$ws.Cells.Item(133, 8).Value = 6783
$ws.Cells.Item(133, 9).Value = 4666.3335
$ws.Cells.Item(133, 10).Value = 8899.666999999999
$ws.Cells.Item(133, 11).Value = 13999.0005
$ws.Cells.Item(133, 12).Value = 26699.001
$ws.Cells.Item(133, 13).Value = -8939.000499999998
$ws.Cells.Item(133, 14).Value = -36819.001
$ws.Cells.Item(134, 8).Value = 6621.531
$ws.Cells.Item(134, 9).Value = 3142.2856
$ws.Cells.Item(134, 10).Value = 7201.405
$ws.Cells.Item(134, 11).Value = 9426.856800000001
$ws.Cells.Item(134, 12).Value = 21604.215
$ws.Cells.Item(134, 13).Value = -4356.856800000001
$ws.Cells.Item(134, 14).Value = -31744.215
$ws.Cells.Item(139, 8).Value = 15444.182
$ws.Cells.Item(139, 9).Value = 17543
$ws.Cells.Item(139, 10).Value = 5999.5
$ws.Cells.Item(139, 11).Value = 52629
$ws.Cells.Item(139, 12).Value = 17998.5
$ws.Cells.Item(139, 13).Value = -47489
$ws.Cells.Item(139, 14).Value = -28278.5
$ws.Cells.Item(141, 8).Value = 20000
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 2810.85
$ws.Cells.Item(46, 8).Value = 1360
$ws.Cells.Item(46, 9).Value = 1360
$ws.Cells.Item(46, 11).Value = 1360
$ws.Cells.Item(46, 13).Value = -1204
$ws.Cells.Item(80, 8).Value = 11382.7
$ws.Cells.Item(80, 10).Value = 10775.23
$ws.Cells.Item(80, 12).Value = 10775.23
$ws.Cells.Item(80, 14).Value = -12771.23
$ws.Cells.Item(83, 8).Value = 11382.7
$ws.Cells.Item(83, 10).Value = 10775.23
$ws.Cells.Item(83, 12).Value = 53876.14999999999
$ws.Cells.Item(83, 14).Value = -63860.14999999999
$ws.Cells.Item(97, 8).Value = 1857.1666
$ws.Cells.Item(97, 9).Value = 1977.5555
$ws.Cells.Item(97, 11).Value = 1977.5555
$ws.Cells.Item(97, 13).Value = -1481.5555
$ws.Cells.Item(126, 8).Value = 9613.522999999999
$ws.Cells.Item(126, 9).Value = 12606
$ws.Cells.Item(126, 10).Value = 6893.091
$ws.Cells.Item(126, 11).Value = 37818
$ws.Cells.Item(126, 12).Value = 20679.273
$ws.Cells.Item(126, 13).Value = -35348
$ws.Cells.Item(126, 14).Value = -25619.273
$ws.Cells.Item(132, 8).Value = 10079.317
$ws.Cells.Item(132, 9).Value = 7856.514
$ws.Cells.Item(132, 10).Value = 23045.666
$ws.Cells.Item(132, 11).Value = 23569.542
$ws.Cells.Item(132, 12).Value = 69136.99800000001
$ws.Cells.Item(132, 13).Value = -21039.542
$ws.Cells.Item(132, 14).Value = -74196.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2439.4324
$ws.Cells.Item(22, 9).Value = 1903.174
$ws.Cells.Item(22, 11).Value = 1903.174
$ws.Cells.Item(22, 13).Value = -1608.174
$ws.Cells.Item(27, 8).Value = 2439.4324
$ws.Cells.Item(27, 9).Value = 1903.174
$ws.Cells.Item(27, 11).Value = 1903.174
$ws.Cells.Item(27, 13).Value = -1796.174
$ws.Cells.Item(55, 8).Value = 2840.5
$ws.Cells.Item(55, 9).Value = 2665.8333
$ws.Cells.Item(55, 10).Value = 2971.5
$ws.Cells.Item(55, 11).Value = 2665.8333
$ws.Cells.Item(55, 12).Value = 2971.5
$ws.Cells.Item(55, 13).Value = -2492.8333
$ws.Cells.Item(55, 14).Value = -3317.5
$ws.Cells.Item(61, 8).Value = 3651.7742
$ws.Cells.Item(61, 9).Value = 2869.0833
$ws.Cells.Item(61, 11).Value = 2869.0833
$ws.Cells.Item(61, 13).Value = -2667.0833
$ws.Cells.Item(68, 8).Value = 1434.8667
$ws.Cells.Item(68, 9).Value = 1139.091
$ws.Cells.Item(68, 10).Value = 2248.25
$ws.Cells.Item(68, 11).Value = 1139.091
$ws.Cells.Item(68, 12).Value = 2248.25
$ws.Cells.Item(68, 13).Value = -390.0909999999999
$ws.Cells.Item(68, 14).Value = -3746.25
$ws.Cells.Item(71, 8).Value = 1434.8667
$ws.Cells.Item(71, 9).Value = 1139.091
$ws.Cells.Item(71, 10).Value = 2248.25
$ws.Cells.Item(71, 11).Value = 5695.455
$ws.Cells.Item(71, 12).Value = 11241.25
$ws.Cells.Item(71, 13).Value = -1951.455
$ws.Cells.Item(71, 14).Value = -18729.25
$ws.Cells.Item(113, 8).Value = 3651.7742
$ws.Cells.Item(113, 9).Value = 2869.0833
$ws.Cells.Item(113, 11).Value = 2869.0833
$ws.Cells.Item(113, 13).Value = -699.0832999999998
$ws.Cells.Item(136, 8).Value = 65519.94
$ws.Cells.Item(136, 9).Value = 119467.82
$ws.Cells.Item(136, 11).Value = 358403.46
$ws.Cells.Item(136, 13).Value = -355853.46

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 17188
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 13).ClearContents()
$ws.Cells.Item(107, 8).Value = 1307.5714
$ws.Cells.Item(107, 9).Value = 1557.3334
$ws.Cells.Item(107, 11).Value = 4672.0002
$ws.Cells.Item(107, 13).Value = -2752.0002
$ws.Cells.Item(122, 8).Value = 3396.4084
$ws.Cells.Item(122, 9).Value = 2605.9333
$ws.Cells.Item(122, 11).Value = 7817.7999
$ws.Cells.Item(122, 13).Value = -5367.7999
$ws.Cells.Item(133, 8).Value = 62748.8
$ws.Cells.Item(133, 10).Value = 62748.8
$ws.Cells.Item(133, 12).Value = 62748.8
$ws.Cells.Item(133, 14).Value = -72868.8
